# Apply "New crime data collected" update to the 76th Precinct weekly CompStat sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text updates (shared-string rich-text runs): bulletin number + dates.
# A8 = "Volume 30   Number  44" -> "...  45"   ("44" run is chars 21-22)
# C9 = "Report Covering the Week  10/30/2023  Through  11/5/2023"
#      -> "...  11/6/2023  Through  11/12/2023"
# ---------------------------------------------------------------------------
$ws.Range("A8").Characters(21, 2).Text = "45"

$c9 = $ws.Range("C9")
$c9.Characters(27, 10).Text = "11/6/2023"
$c9.Characters(47, 9).Text = "11/12/2023"

# ---------------------------------------------------------------------------
# Stable "donor" cells used purely to copy a number-format/style onto a cell
# whose value is switching between text-placeholder ("0" / "***.*") and a
# real number (or vice versa). None of these donor cells are themselves
# edited below.
#   D17 -> style for plain integer counts      (#,##0)
#   E18 -> style for percentage-change figures (#,##0.0;"-"#,##0.0)
#   C26 -> style for the literal text "0"
#   E26 -> style for the literal text "***.*"
# ---------------------------------------------------------------------------
$intStyle = $ws.Range("D17")
$pctStyle = $ws.Range("E18")
$zeroText = $ws.Range("C26")
$starText = $ws.Range("E26")

# Row 16 - Robbery: was all text placeholders in C/D/E, now real numbers.
$intStyle.Copy($ws.Range("C16"))
$ws.Range("C16").Value = 2
$intStyle.Copy($ws.Range("D16"))
$ws.Range("D16").Value = 1
$pctStyle.Copy($ws.Range("E16"))
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 3
$ws.Range("G16").Value = 2
$ws.Range("H16").Value = 50
$ws.Range("I16").Value = 57
$ws.Range("J16").Value = 51
$ws.Range("K16").Value = 11.764705882352
$ws.Range("L16").Value = 83.870967741935
$ws.Range("M16").Value = -38.709677419354
$ws.Range("N16").Value = -83.136094674556

# Row 17 - Felonious Assault
$ws.Range("C17").Value = 2
$ws.Range("E17").Value = 0
$ws.Range("G17").Value = 14
$ws.Range("H17").Value = -57.142857142857
$ws.Range("I17").Value = 91
$ws.Range("J17").Value = 90
$ws.Range("K17").Value = 1.111111111111
$ws.Range("L17").Value = 13.75
$ws.Range("M17").Value = 12.345679012345
$ws.Range("N17").Value = -66.296296296296

# Row 18 - Burglary
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 4
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = -70
$ws.Range("I18").Value = 79
$ws.Range("J18").Value = 95
$ws.Range("K18").Value = -16.842105263157
$ws.Range("L18").Value = 27.419354838709
$ws.Range("M18").Value = -13.186813186813
$ws.Range("N18").Value = -78.877005347593

# Row 19 - Grand Larceny
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 3
$ws.Range("E19").Value = 33.333333333333
$ws.Range("F19").Value = 13
$ws.Range("G19").Value = 22
$ws.Range("H19").Value = -40.909090909090
$ws.Range("I19").Value = 142
$ws.Range("J19").Value = 185
$ws.Range("K19").Value = -23.243243243243
$ws.Range("L19").Value = 25.663716814159
$ws.Range("M19").Value = -27.551020408163
$ws.Range("N19").Value = -32.380952380952

# Row 20 - G.L.A.
$ws.Range("C20").Value = 2
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 6
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 56
$ws.Range("J20").Value = 52
$ws.Range("K20").Value = 7.692307692307
$ws.Range("L20").Value = 14.285714285714
$ws.Range("M20").Value = 1.818181818181
$ws.Range("N20").Value = -84.135977337110

# Row 21 - TOTAL
$ws.Range("C21").Value = 12
$ws.Range("D21").Value = 12
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 31
$ws.Range("G21").Value = 55
$ws.Range("H21").Value = -43.636363636363
$ws.Range("I21").Value = 431
$ws.Range("J21").Value = 480
$ws.Range("K21").Value = -10.208333333333
$ws.Range("L21").Value = 26.023391812865
$ws.Range("M21").Value = -16.955684007707
$ws.Range("N21").Value = -72.442455242966

# Row 22 - Transit (only L22 changes)
$ws.Range("L22").Value = -71.428571428571

# Row 23 - Housing: C goes text -> number, D & E go number -> text.
$intStyle.Copy($ws.Range("C23"))
$ws.Range("C23").Value = 5
$zeroText.Copy($ws.Range("D23"))
$starText.Copy($ws.Range("E23"))
$ws.Range("F23").Value = 5
$ws.Range("H23").Value = 25
$ws.Range("I23").Value = 93
$ws.Range("K23").Value = -4.123711340206
$ws.Range("L23").Value = 22.368421052631
$ws.Range("M23").Value = 30.985915492957

# Row 24 - Petit Larceny
$ws.Range("C24").Value = 21
$ws.Range("E24").Value = 162.5
$ws.Range("F24").Value = 49
$ws.Range("G24").Value = 33
$ws.Range("H24").Value = 48.484848484848
$ws.Range("I24").Value = 536
$ws.Range("J24").Value = 396
$ws.Range("K24").Value = 35.353535353535
$ws.Range("L24").Value = 73.462783171521
$ws.Range("M24").Value = 35.353535353535

# Row 25 - Misdemeanor Assault
$ws.Range("F25").Value = 11
$ws.Range("G25").Value = 4
$ws.Range("H25").Value = 175
$ws.Range("I25").Value = 138
$ws.Range("J25").Value = 151
$ws.Range("K25").Value = -8.609271523178
$ws.Range("L25").Value = -1.428571428571
$ws.Range("M25").Value = -45.669291338582

# Row 27 - Other Sex Crimes: D/E go text -> number.
$intStyle.Copy($ws.Range("D27"))
$ws.Range("D27").Value = 1
$pctStyle.Copy($ws.Range("E27"))
$ws.Range("E27").Value = -100
$ws.Range("J27").Value = 18
$ws.Range("K27").Value = -38.888888888888

# Row 28 - Shooting Victims: C/F go text -> number.
$intStyle.Copy($ws.Range("C28"))
$ws.Range("C28").Value = 1
$intStyle.Copy($ws.Range("F28"))
$ws.Range("F28").Value = 1
$ws.Range("I28").Value = 8
$ws.Range("K28").Value = 33.333333333333
$ws.Range("L28").Value = -33.333333333333
$ws.Range("M28").Value = -50
$ws.Range("N28").Value = -85.185185185185

# Row 29 - Shooting Incidents: C/F go text -> number.
$intStyle.Copy($ws.Range("C29"))
$ws.Range("C29").Value = 1
$intStyle.Copy($ws.Range("F29"))
$ws.Range("F29").Value = 1
$ws.Range("I29").Value = 8
$ws.Range("K29").Value = 100
$ws.Range("L29").Value = -27.272727272727
$ws.Range("M29").Value = -38.461538461538
$ws.Range("N29").Value = -81.395348837209
